# Add 2022-Q3 data:
#  1. Insert a new "2022-Q3" row at the top of the "总计" (summary) sheet,
#     pushing the existing quarters down by one row.
#  2. Insert a brand-new "2022-Q3" worksheet (with the per-fund holding
#     detail) right after "总计" and before the existing "2022-Q2" sheet.

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)

# --- 1. Update the "总计" summary sheet -----------------------------------
$wsTotal.Rows(2).Insert()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 5
$wsTotal.Range("D2").Value = 0.2

# Restore the formatting that Insert() did not carry onto the new row:
#  - A2 should use the bold/bordered "index" style used by the other rows
#  - B2:D2 should have no special formatting (Insert() left them styled)
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$wsTotal.Range("B2:D2").ClearFormats()

# --- 2. Insert the new "2022-Q3" detail worksheet --------------------------
# Worksheets.Add(Before) inserts immediately before the given sheet, which
# places the new sheet right after "总计" (i.e. before the current "2022-Q2").
$wsQ3 = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$wsQ3.Name = "2022-Q3"

$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

$data = @(
  @(0, "012586", "南方港股创新视野一年持有混合A", "2.16", "50.74", "3.56", "0.0769", 5),
  @(1, "470888", "汇添富香港优势精选混合（QDII）", "1.63", "78.50", "4.32", "0.0704", 7),
  @(2, "004266", "招商沪港深科技创新主题精选灵活配置混合A", "0.92", "90.52", "4.17", "0.0384", 2),
  @(3, "010754", "招商沪港深科技创新主题精选灵活配置混合C", "0.25", "90.52", "4.17", "0.0104", 2),
  @(4, "012587", "南方港股创新视野一年持有混合C", "0.19", "50.74", "3.56", "0.0068", 5)
)

$r = 2
foreach ($row in $data) {
  $wsQ3.Range("A$r").Value = $row[0]
  $wsQ3.Range("B$r").Value = "'" + $row[1]
  $wsQ3.Range("C$r").Value = $row[2]
  $wsQ3.Range("D$r").Value = "'" + $row[3]
  $wsQ3.Range("E$r").Value = "'" + $row[4]
  $wsQ3.Range("F$r").Value = "'" + $row[5]
  $wsQ3.Range("G$r").Value = "'" + $row[6]
  $wsQ3.Range("H$r").Value = $row[7]
  $r = $r + 1
}

# Match the bold/bordered header + index-column styling used on every other
# quarter sheet by copying it from the "总计" sheet's equivalent cells.
$wsTotal.Range("B1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A6").PasteSpecial(-4122)
